$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts replacing the previous Strike# derived values.
$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 2
    24 = 3
    25 = 1
    26 = 2
    27 = 1
    28 = 2
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 1
    38 = 0
    39 = 0
    40 = 0
    41 = 2
    42 = 1
    43 = 0
    44 = 0
    45 = 0
    46 = 1
    47 = 2
    48 = 1
    49 = 0
    50 = 2
    51 = 2
    52 = 1
    54 = 2
    55 = 0
    56 = 0
    57 = 1
    58 = 2
    59 = 1
    60 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
